{"js": "// Load all paragraphs in the body so we can address them by index.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Remove the first paragraph entirely ---\n// (\"W datetime picker wy\u015bwietlaj data dzisiejsza podczas dodawania nowego wniosku\")\nparagraphs.items[0].delete();\nawait context.sync();\n\n// Re-fetch paragraphs after the deletion shifted indices.\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 2. Rewrite the \"Object reference...\" paragraph (now index 0) ---\n// Collapse its many spell-checked runs into a single run, then append the\n// extra \" co to \" / \"dok\u0142adnie \" / \"znaczy\" runs that replace the old\n// trailing \" co to znaczy\" text.\nconst objParagraph = paragraphs.items[0];\nobjParagraph.insertText(\n  \"Object reference not set to an instance of an object\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\nobjParagraph.insertText(\" co to \", Word.InsertLocation.end);\nawait context.sync();\nobjParagraph.insertText(\"dok\u0142adnie \", Word.InsertLocation.end);\nawait context.sync();\nobjParagraph.insertText(\"znaczy\", Word.InsertLocation.end);\nawait context.sync();\n\n// --- 3. Collapse the \"Doda\u0107 status wniosku...\" paragraph (now index 1) ---\n// into a single plain run (same text, just without the spell-check run\n// splits / proofErr markers).\nparagraphs.load(\"items\");\nawait context.sync();\nconst statusParagraph = paragraphs.items[1];\nstatusParagraph.insertText(\n  \"Doda\u0107 status wniosku: Z\u0142o\u017cony, Oferta, Decyzja, Wyp\u0142acony, Negat, Rezygnacja Klienta\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 4. \"Doda\u0107 Status klienta...\" (now index 2) is left untouched. ---\n\n// --- 5. Insert the new task paragraphs after the \"Status klienta\" paragraph\n// and before the trailing empty paragraph. ---\nparagraphs.load(\"items\");\nawait context.sync();\nconst clientStatusParagraph = paragraphs.items[2];\n\nconst newParagraphTexts = [\n  \"Zaprogramowa\u0107 przycisk dodaj nowy miesi\u0105c\",\n  \"Doda\u0107 list\u0119 po\u015brednik\u00f3w (podobna lista do ZusUsList)\",\n  \"Doda\u0107 list\u0119 multibroker\u00f3w do Grida Result albo doda\u0107 list\u0119 do Wniosku\",\n  \"Zaprogramowa\u0107 dodawanie zada\u0144 do klient\u00f3w i wniosk\u00f3w\",\n  \"Podzieli\u0107 wnioski wy\u015bwietlane w LoanSearch na miesi\u0105ce\",\n  \"Zaprogramuj wszystkie kontrolki w ClientDetails\",\n];\n\nlet anchor = clientStatusParagraph;\nfor (const text of newParagraphTexts) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Last new paragraph has two runs: \"Dodaj przycisk szczeg\u00f3\u0142y do Grid Loan w \"\n// + \"ClientDetails\".\nconst lastParagraph = anchor.insertParagraph(\n  \"Dodaj przycisk szczeg\u00f3\u0142y do Grid Loan w \",\n  Word.InsertLocation.after\n);\nawait context.sync();\nlastParagraph.insertText(\"ClientDetails\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Remove the first paragraph entirely ---\n# (\"W datetime picker wy\u015bwietlaj data dzisiejsza podczas dodawania nowego wniosku\")\n$d.Paragraphs(1).Range.Delete()\n\n# --- 2. Rewrite the \"Object reference...\" paragraph (now paragraph 1) ---\n# Collapse its many spell-checked runs into a single run, then append the\n# extra \" co to \" / \"dok\u0142adnie \" / \"znaczy\" runs that replace the old\n# trailing \" co to znaczy\" text.\n$objRange = $d.Paragraphs(1).Range\n$objRange.Find.Execute(\n    \"Object reference not set to an instance of an object co to znaczy\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Object reference not set to an instance of an object\", 2) | Out-Null\n\n$objEnd = $d.Paragraphs(1).Range\n$ip1 = $d.Range($objEnd.End - 1, $objEnd.End - 1)\n$ip1.InsertAfter(\" co to \")\n\n$objEnd2 = $d.Paragraphs(1).Range\n$ip2 = $d.Range($objEnd2.End - 1, $objEnd2.End - 1)\n$ip2.InsertAfter(\"dok\u0142adnie \")\n\n$objEnd3 = $d.Paragraphs(1).Range\n$ip3 = $d.Range($objEnd3.End - 1, $objEnd3.End - 1)\n$ip3.InsertAfter(\"znaczy\")\n\n# --- 3. Collapse the \"Doda\u0107 status wniosku...\" paragraph (now paragraph 2) ---\n# into a single plain run (same text, just without the spell-check run\n# splits / proofErr markers). Use Find/Replace (rather than a plain Range.Text\n# assignment) so the whole multi-run paragraph range is replaced instead of\n# only its first run.\n$statusRange = $d.Paragraphs(2).Range\n$statusRange.Find.Execute(\n    \"Doda\u0107 status wniosku: Z\u0142o\u017cony, Oferta, Decyzja, Wyp\u0142acony, Negat, Rezygnacja Klienta\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Doda\u0107 status wniosku: Z\u0142o\u017cony, Oferta, Decyzja, Wyp\u0142acony, Negat, Rezygnacja Klienta\", 2) | Out-Null\n\n# --- 4. \"Doda\u0107 Status klienta...\" (paragraph 3) is left untouched. ---\n\n# --- 5. Insert the new task paragraphs after the \"Status klienta\" paragraph\n# and before the trailing empty paragraph. ---\n$newParagraphTexts = @(\n    \"Zaprogramowa\u0107 przycisk dodaj nowy miesi\u0105c\",\n    \"Doda\u0107 list\u0119 po\u015brednik\u00f3w (podobna lista do ZusUsList)\",\n    \"Doda\u0107 list\u0119 multibroker\u00f3w do Grida Result albo doda\u0107 list\u0119 do Wniosku\",\n    \"Zaprogramowa\u0107 dodawanie zada\u0144 do klient\u00f3w i wniosk\u00f3w\",\n    \"Podzieli\u0107 wnioski wy\u015bwietlane w LoanSearch na miesi\u0105ce\",\n    \"Zaprogramuj wszystkie kontrolki w ClientDetails\"\n)\n\n$insertAfterIndex = 3\nforeach ($text in $newParagraphTexts) {\n    $d.Paragraphs($insertAfterIndex).Range.InsertParagraphAfter()\n    $insertAfterIndex = $insertAfterIndex + 1\n    $d.Paragraphs($insertAfterIndex).Range.Text = $text\n}\n\n# Last new paragraph has two runs: \"Dodaj przycisk szczeg\u00f3\u0142y do Grid Loan w \"\n# + \"ClientDetails\".\n$d.Paragraphs($insertAfterIndex).Range.InsertParagraphAfter()\n$insertAfterIndex = $insertAfterIndex + 1\n$d.Paragraphs($insertAfterIndex).Range.Text = \"Dodaj przycisk szczeg\u00f3\u0142y do Grid Loan w \"\n$lastRange = $d.Paragraphs($insertAfterIndex).Range\n$lastIp = $d.Range($lastRange.End - 1, $lastRange.End - 1)\n$lastIp.InsertAfter(\"ClientDetails\")\n"}
